$d = $word.ActiveDocument

# --- Table 1 (the "a) na przepustkę jednorazową" table) ---
$t1 = $d.Tables.Item(1)

# Remove the old row "1)" (sierż. pchor. / Paweł / BYŚ / 26.05.2021 / Kraków)
$t1.Rows.Item(1).Delete()

# Renumber the remaining rows: old "2)" -> "1)", old "3)" -> "2)"
$t1.Rows.Item(1).Cells.Item(1).Range.Text = "1)"
$t1.Rows.Item(2).Cells.Item(1).Range.Text = "2)"

# Fix up the corrupted city name for the row that is now "1)"
$t1.Rows.Item(1).Cells.Item(7).Range.Text = "Warszawhgfghva"

# Append a brand-new row "3)" with a new travelling cadet entry
$newRow1 = $t1.Rows.Add()
$newRow1.Cells.Item(1).Range.Text = "3)"
$newRow1.Cells.Item(2).Range.Text = "szer. pchor."
$newRow1.Cells.Item(3).Range.Text = "Paweł"
$newRow1.Cells.Item(4).Range.Text = "BYŚ"
$newRow1.Cells.Item(5).Range.Text = "w dn. 05 - 07.11.2021 r."
$newRow1.Cells.Item(6).Range.Text = "do m."
$newRow1.Cells.Item(7).Range.Text = "Dębica"

# --- Table 2 (the "b) na urlop" table) ---
$t2 = $d.Tables.Item(2)

# Append a brand-new row "5)" with a new travelling cadet entry
$newRow2 = $t2.Rows.Add()
$newRow2.Cells.Item(1).Range.Text = "5)"
$newRow2.Cells.Item(2).Range.Text = "szer. pchor."
$newRow2.Cells.Item(3).Range.Text = "Paweł"
$newRow2.Cells.Item(4).Range.Text = "BYŚ"
$newRow2.Cells.Item(5).Range.Text = "w dn. 06 - 08.11.2021 r."
$newRow2.Cells.Item(6).Range.Text = "do m."
$newRow2.Cells.Item(7).Range.Text = "Kraków"
